# Update TPM-derived NATMI ligand-receptor metrics (Vwf-Tnfrsf11b) with
# refreshed values produced by the updated TPM computation scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 89.35833500000001
$ws.Range("H2").Value = 268.075005
$ws.Range("I2").Value = 0.9624640326757887
$ws.Range("J2").Value = 0.9624640326757889
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("Q2").Value = 7.364228890131668
$ws.Range("R2").Value = 66.278060011185
$ws.Range("S2").Value = 0.04549572884635764
$ws.Range("T2").Value = 0.04549572884635765
$ws.Range("G3").Value = 89.35833500000001
$ws.Range("H3").Value = 268.075005
$ws.Range("I3").Value = 0.9624640326757887
$ws.Range("J3").Value = 0.9624640326757889
$ws.Range("Q3").Value = 148.42633903504
$ws.Range("R3").Value = 1335.83705131536
$ws.Range("S3").Value = 0.9169683038294312
$ws.Range("T3").Value = 0.9169683038294313
$ws.Range("I4").Value = 0.001854741667334279
$ws.Range("J4").Value = 0.001854741667334279
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("S4").Value = 0.00008767374271897227
$ws.Range("T4").Value = 0.00008767374271897229
$ws.Range("I5").Value = 0.001854741667334279
$ws.Range("J5").Value = 0.001854741667334279
$ws.Range("S5").Value = 0.001767067924615307
$ws.Range("T5").Value = 0.001767067924615307
$ws.Range("G6").Value = 1.963978
$ws.Range("H6").Value = 5.891934
$ws.Range("I6").Value = 0.02115368628977398
$ws.Range("J6").Value = 0.02115368628977398
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("Q6").Value = 0.1618560095953333
$ws.Range("R6").Value = 1.456704086358
$ws.Range("S6").Value = 0.0009999359382447286
$ws.Range("T6").Value = 0.0009999359382447288
$ws.Range("G7").Value = 1.963978
$ws.Range("H7").Value = 5.891934
$ws.Range("I7").Value = 0.02115368628977398
$ws.Range("J7").Value = 0.02115368628977398
$ws.Range("Q7").Value = 3.262214593472
$ws.Range("R7").Value = 29.359931341248
$ws.Range("S7").Value = 0.02015375035152925
$ws.Range("T7").Value = 0.02015375035152925
$ws.Range("G8").Value = 0.16825
$ws.Range("H8").Value = 0.50475
$ws.Range("I8").Value = 0.001812193272151965
$ws.Range("J8").Value = 0.001812193272151965
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("Q8").Value = 0.01386587508333333
$ws.Range("R8").Value = 0.12479287575
$ws.Range("S8").Value = 0.00008566247769052179
$ws.Range("T8").Value = 0.0000856624776905218
$ws.Range("G9").Value = 0.16825
$ws.Range("H9").Value = 0.50475
$ws.Range("I9").Value = 0.001812193272151965
$ws.Range("J9").Value = 0.001812193272151965
$ws.Range("Q9").Value = 0.279467288
$ws.Range("R9").Value = 2.515205592
$ws.Range("S9").Value = 0.001726530794461443
$ws.Range("T9").Value = 0.001726530794461443
$ws.Range("G10").Value = 1.180534666666667
$ws.Range("H10").Value = 3.541604
$ws.Range("I10").Value = 0.01271534609495094
$ws.Range("J10").Value = 0.01271534609495094
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("Q10").Value = 0.09729061646088889
$ws.Range("R10").Value = 0.875615548148
$ws.Range("S10").Value = 0.000601055123603096
$ws.Range("T10").Value = 0.0006010551236030962
$ws.Range("G11").Value = 1.180534666666667
$ws.Range("H11").Value = 3.541604
$ws.Range("I11").Value = 0.01271534609495094
$ws.Range("J11").Value = 0.01271534609495094
$ws.Range("Q11").Value = 1.960896414165334
$ws.Range("R11").Value = 17.648067727488
$ws.Range("S11").Value = 0.01211429097134785
$ws.Range("T11").Value = 0.01211429097134785
